$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new "Prisma_Population" column is inserted in front of the existing
# data, pushing every existing column one letter to the right (old A->B,
# old B->C, ... old G->H). Using a real column insert (rather than
# clearing + retyping every cell) preserves the untouched columns' widths
# and styles exactly, and naturally keeps the Study_Types column (old B,
# now C) lined up with its original per-row value.
$ws.Columns.Item(1).Insert()

# --- New column A: population name for each block of 4 rows ---
$ws.Range("A1").Value = "Prisma_Population"

# --- Introduce the ICER-block shared strings first (rows 7-10), matching
#     the original authoring order where this block's distinct text was
#     entered before the Test_Sachin block's distinct text. ---
$ws.Range("A7").Value = "ICER RRMM 2022 report"
$ws.Range("B7").Value = "\Testdata\Templates\PRISMA\ICER\5. PRISMA.xlsx"
$ws.Range("H7").Value = "\Testdata\Templates\PRISMA\ICER\5. PRISMA_ICER.Clinical.PNG"
$ws.Range("H9").Value = "\Testdata\Templates\PRISMA\ICER\5. PRISMA_ICER.QOL.PNG"
$ws.Range("H8").Value = "\Testdata\Templates\PRISMA\ICER\5. PRISMA_ICER.Econ.PNG"
$ws.Range("H10").Value = "\Testdata\Templates\PRISMA\ICER\5. PRISMA_ICER.Clinical.PNG"
$ws.Range("A8").Value = "ICER RRMM 2022 report"
$ws.Range("A9").Value = "ICER RRMM 2022 report"
$ws.Range("A10").Value = "ICER RRMM 2022 report"

# --- Introduce the Test_Sachin-block shared strings next (rows 2-5) ---
$ws.Range("A2").Value = "Test_Sachin"
$ws.Range("B2").Value = "\Testdata\Templates\PRISMA\Test_Sachin\13. PRISMA_AAA_mCRPC_ID Update.xlsx"
$ws.Range("H2").Value = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Clinical.PNG"
$ws.Range("H3").Value = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Econ.PNG"
$ws.Range("H4").Value = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_QoL.PNG"
$ws.Range("H5").Value = "\Testdata\Templates\PRISMA\Test_Sachin\mCRPC_Clinical.PNG"
$ws.Range("A3").Value = "Test_Sachin"
$ws.Range("A4").Value = "Test_Sachin"
$ws.Range("A5").Value = "Test_Sachin"

# --- Remove the stale old column-B cells that landed on rows 3-5 after
#     the insert (old row 3 had a Study_Types cell in column B; rows 4-5
#     never had one, so nothing to remove there). The new layout has no
#     column-B entry on rows 3-5 (Study_Types now lives in column C and
#     column B is only populated on rows 2 and 7). ---
$ws.Range("B3").Clear()

# --- New numeric columns D:G for the Test_Sachin block ---
$ws.Range("D2").Value = 500
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 1500
$ws.Range("G2").Value = 2000

$ws.Range("D3").Value = 600
$ws.Range("E3").Value = 1200
$ws.Range("F3").Value = 1800
$ws.Range("G3").Value = 2400

$ws.Range("D4").Value = 700
$ws.Range("E4").Value = 1400
$ws.Range("F4").Value = 2100
$ws.Range("G4").Value = 2800

$ws.Range("D5").Value = 800
$ws.Range("E5").Value = 1600
$ws.Range("F5").Value = 2400
$ws.Range("G5").Value = 3200

# --- Numeric + Study_Types columns for the new ICER block (rows 7-10) ---
$ws.Range("C7").Value = "Interventional"
$ws.Range("D7").Value = 501
$ws.Range("E7").Value = 1001
$ws.Range("F7").Value = 1501
$ws.Range("G7").Value = 2001

$ws.Range("C8").Value = "Economic"
$ws.Range("D8").Value = 601
$ws.Range("E8").Value = 1201
$ws.Range("F8").Value = 1801
$ws.Range("G8").Value = 2401

$ws.Range("C9").Value = "Quality of life"
$ws.Range("D9").Value = 701
$ws.Range("E9").Value = 1401
$ws.Range("F9").Value = 2101
$ws.Range("G9").Value = 2801

$ws.Range("C10").Value = "Real-world Evidence"
$ws.Range("D10").Value = 801
$ws.Range("E10").Value = 1601
$ws.Range("F10").Value = 2401
$ws.Range("G10").Value = 3201

# --- Column A (new) uses the sheet's default width, same as before the
#     insert nudged everything -- give it no custom width. Columns B and H
#     need to widen to fit their new (longer) file-path contents. ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(8).EntireColumn.AutoFit()

# --- View: scrolled so column C is leftmost, with H4 selected ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H4").Select()
